$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 368, shifting existing rows 368..455 down to 369..456
$ws.Rows.Item(368).Insert()

# Populate the newly inserted row 368 with the new data record
$ws.Cells.Item(368, 1).Value = 9
$ws.Cells.Item(368, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(368, 3).Value = "Metropolitana"
$ws.Cells.Item(368, 4).Value = 45211
$ws.Cells.Item(368, 5).Value = 13
$ws.Cells.Item(368, 6).Value = 100112043
$ws.Cells.Item(368, 7).Value = "Pepino ensalada"
$ws.Cells.Item(368, 8).Value = "Sin especificar"
$ws.Cells.Item(368, 9).Value = "Primera"
$ws.Cells.Item(368, 10).Value = 70
$ws.Cells.Item(368, 11).Value = 14000
$ws.Cells.Item(368, 12).Value = 16000
$ws.Cells.Item(368, 13).Value = 15029
$ws.Cells.Item(368, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(368, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(368, 16).Value = 250
$ws.Cells.Item(368, 17).Value = 60
$ws.Cells.Item(368, 18).Value = "Hortaliza"
